$d = $word.ActiveDocument

$replacements = @(
    @{old="731×6="; new="660×8="},
    @{old="891×6="; new="375×3="},
    @{old="691×6="; new="931×6="},
    @{old="271×3="; new="913×5="},
    @{old="514×9="; new="565×9="},
    @{old="134×9="; new="839×2="},
    @{old="694×9="; new="368×8="},
    @{old="496×5="; new="588×5="},
    @{old="573×2="; new="816×8="},
    @{old="150×6="; new="507×4="},
    @{old="915×8="; new="869×6="},
    @{old="952×5="; new="807×3="},
    @{old="244×2="; new="215×7="},
    @{old="182×5="; new="832×4="},
    @{old="951×6="; new="316×9="},
    @{old="876×8="; new="817×6="},
    @{old="550×6="; new="772×5="},
    @{old="299×2="; new="549×9="},
    @{old="902×6="; new="748×6="},
    @{old="882×7="; new="288×4="},
    @{old="752×5="; new="341×8="},
    @{old="844×6="; new="977×9="},
    @{old="976×5="; new="254×5="},
    @{old="410×6="; new="548×5="},
    @{old="888×3="; new="679×6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
